# Insert a new data row at row 128, shifting the existing rows 128-237 down
# to 129-238, then populate the newly inserted row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 128 (pushes 128..237 -> 129..238)
$ws.Rows(128).Insert()

# Fill in the new row 128 with the new record's data.
$ws.Cells.Item(128, 1).Value = 9
$ws.Cells.Item(128, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(128, 3).Value = "Metropolitana"
$ws.Cells.Item(128, 4).Value = 44596
$ws.Cells.Item(128, 5).Value = 13
$ws.Cells.Item(128, 6).Value = 100112001
$ws.Cells.Item(128, 7).Value = "Berenjena"
$ws.Cells.Item(128, 8).Value = "Sin especificar"
$ws.Cells.Item(128, 9).Value = "Primera"
$ws.Cells.Item(128, 10).Value = 130
$ws.Cells.Item(128, 11).Value = 12000
$ws.Cells.Item(128, 12).Value = 12000
$ws.Cells.Item(128, 13).Value = 12000
$ws.Cells.Item(128, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(128, 15).Value = "Región Metropolitana"
$ws.Cells.Item(128, 16).Value = 200
$ws.Cells.Item(128, 17).Value = 60
$ws.Cells.Item(128, 18).Value = "Hortaliza"
